$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Cells.Item(1, 19).Value = "IF_val"
$ws.Cells.Item(1, 20).Value = "IF_chg"

$special = @{
    7 = @(0.04, -0.27999999999999997, $true)
    10 = @(-0.15, 0.07999999999999999, $false)
    14 = @(0.02, -0.06, $true)
    27 = @(0.02, -0.16999999999999998, $true)
    30 = @(-0.01, 0.25, $true)
}

for ($r = 2; $r -le 35; $r++) {
    if ($special.ContainsKey($r)) {
        $vals = $special[$r]
        $sVal = $vals[0]
        $tVal = $vals[1]
        $leftAlign = $vals[2]
        $sCell = $ws.Cells.Item($r, 19)
        $sCell.Value = $sVal
        $sCell.NumberFormat = "0.00"
        if ($leftAlign) {
            $sCell.HorizontalAlignment = -4131
        }
        $ws.Cells.Item($r, 20).Value = $tVal
    } else {
        $ws.Cells.Item($r, 19).Value = 0
        $ws.Cells.Item($r, 19).NumberFormat = "0.00"
        $ws.Cells.Item($r, 20).Value = 0
        $ws.Cells.Item($r, 20).NumberFormat = "0.00"
    }
}

# Row 36: empty formatted cells
$ws.Cells.Item(36, 19).NumberFormat = "0.00"
$ws.Cells.Item(36, 20).NumberFormat = "0.00"

# Selection
$ws.Range("V21").Select()
